$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo" shifts from D to E)
$ws.Range("D1").EntireColumn.Insert()

# Copy the style of the header row (e.g. C1) to the new header cell D1
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# Set the new header and value
$ws.Range("D1").Value = "MAE"
$ws.Range("D2").Value = 1.350102492671558
